$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Name" -> "Full Name"
$ws.Range("A1").Value = "Full Name"

# Update selection to A2
$ws.Range("A2").Select()
